$d = $word.ActiveDocument

$pairs = @(
    @("2024-08-28 Wednesday", "2024-08-29 Thursday"),
    @("42÷2=", "60÷8="),
    @("12÷3=", "13÷8="),
    @("60÷2=", "38÷5="),
    @("27÷7=", "64÷4="),
    @("18÷8=", "39÷6="),
    @("16÷2=", "20÷7="),
    @("69÷4=", "70÷8="),
    @("45÷6=", "74÷7="),
    @("35÷3=", "72÷5="),
    @("90÷2=", "19÷5="),
    @("87÷7=", "84÷8="),
    @("95÷3=", "29÷5="),
    @("23÷5=", "74÷8="),
    @("36÷4=", "39÷6="),
    @("91÷4=", "59÷3="),
    @("54÷8=", "97÷2="),
    @("29÷2=", "30÷6="),
    @("70÷6=", "33÷4="),
    @("22÷8=", "74÷2="),
    @("80÷7=", "28÷2="),
    @("20÷2=", "54÷7="),
    @("85÷5=", "47÷5="),
    @("10÷6=", "28÷2="),
    @("20÷4=", "80÷6="),
    @("25÷5=", "87÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
